$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header and row labels
$ws.Range("B1").Value = "Price per ft^2"
$ws.Range("A2").Value = "San Francisco, CA"
$ws.Range("A3").Value = "Belhaven, NC"
$ws.Range("A4").Value = "Dallas, TX"
$ws.Range("A5").Value = "New York City, NY"

# Update values
$ws.Range("B2").Value = 728
$ws.Range("B3").Value = 28
$ws.Range("B4").Value = 184
$ws.Range("B5").Value = 650

# Widen column A (COM ColumnWidth is stored in character units and gets
# quantized to the Normal-style font's pixel grid on save, same as real
# Excel; 19.83 is the closest input to the target stored width of
# 20.6640625 twips/256).
$ws.Columns.Item(1).ColumnWidth = 19.83

# Update selection to B5
$ws.Range("B5").Select()
